$d = $word.ActiveDocument

# ------------------------------------------------------------------
# This document has a paragraph that ends with "...head over to the
# next lesson where " + a "_GoBack" bookmark (collapsed) + the run
# "we're going to get started building out our Flash Chat app.",
# immediately followed by a paragraph that only contains a legacy
# ActiveX/OLE control object, and finally a trailing empty paragraph.
#
# The edit:
#   1. moves the "_GoBack" bookmark so it sits AFTER the
#      "...Flash Chat app." sentence, in its own (new) paragraph;
#   2. removes the paragraph that hosts the ActiveX/OLE control.
# ------------------------------------------------------------------

# Locate the paragraph that hosts the legacy ActiveX/OLE control.
# It renders as a single (placeholder) character and carries a
# distinctive shading fill (F2F3F5) that is different from the
# body-text shading (F7F8FA) and from the plain trailing paragraph
# (no fill at all) - that combination uniquely identifies it.
$olePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Length -eq 1 -and $para.Range.Shading.BackgroundPatternColor -eq 16118770) {
        $olePara = $para
    }
}

# Step 1: pull the "_GoBack" bookmark out of its current location
# (right before "we're going to get started...").
$hadGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hadGoBack) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

# Step 2: delete the whole ActiveX/OLE-control paragraph. This merges
# it away, so the paragraph that used to follow it (the plain,
# trailing empty paragraph) becomes the new last paragraph of the
# document body.
if ($olePara -ne $null) {
    $olePara.Range.Delete()
}

# Step 3: re-create the "_GoBack" bookmark, collapsed, at the start
# of the (now) last paragraph, so it lands in its own paragraph right
# after "...Flash Chat app." Temporarily typing a character into that
# paragraph first and removing it afterwards keeps the bookmark
# anchored correctly once the placeholder text is gone.
if ($hadGoBack) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Range.InsertAfter("x")
    $anchor = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
    $d.Bookmarks.Add("_GoBack", $anchor)
    $placeholder = $d.Range($lastPara.Range.Start, $lastPara.Range.Start + 1)
    $placeholder.Delete()
}
